$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.267.66"
$ws.Range("E2").Value = "  +5.86%  "
$ws.Range("D3").Value = "3.548.93"
$ws.Range("E3").Value = "  +5.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.21"
$ws.Range("E5").Value = "  +9.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "564.44"
$ws.Range("E6").Value = "  +6.34%  "
$ws.Range("D7").Value = "3.543.94"
$ws.Range("E7").Value = "  +5.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.619"
$ws.Range("E8").Value = "  +3.82%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.633"
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("E11").Value = "  +14.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.80"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +6.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.36"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "4.113.09"
$ws.Range("E15").Value = "  +6.18%  "
$ws.Range("D16").Value = "3.558.69"
$ws.Range("E16").Value = "  +6.50%  "
$ws.Range("D19").Value = "67.282.67"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.06"
$ws.Range("E20").Value = "  +7.50%  "
$ws.Range("E21").Value = "  +3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "424.95"
$ws.Range("E22").Value = "  +14.12%  "
$ws.Range("E23").Value = "  +10.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.42"
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.10"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.91"
$ws.Range("E27").Value = "  +7.52%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.26"
$ws.Range("E29").Value = "  +8.39%  "
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.53"
$ws.Range("E31").Value = "  +5.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "632.01"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.74"
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").Value = "  +5.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.32"
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.42"
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("E38").Value = "  +18.96%  "
$ws.Range("E39").Value = "  +12.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0418"
$ws.Range("E48").Value = "  +4.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +5.47%  "
# Row 17/18 swap (Chainlink <-> TRON)
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.122"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.54"
$ws.Range("E18").Value = "  +5.82%  "

# Row 42/43 swap (Stacks <-> Maker)
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.154.59"
$ws.Range("E42").Value = "  +7.65%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("E43").Value = "  +11.85%  "

# Row 46/47 swap (ThetaToken <-> ApeXProtocol)
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.42"
$ws.Range("E46").Value = "  +12.12%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.87"
$ws.Range("E47").Value = "  +9.94%  "

# Row 51 full replacement (THORChain -> Monero)
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.37"
$ws.Range("E51").Value = "  +3.77%  "
